$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '67.891.79'
Set-TextValue 'E2' '  -1.54%  '
Set-TextValue 'D3' '3.271.25'
Set-TextValue 'E3' '  -1.25%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '185.95'
Set-TextValue 'E5' '  -0.68%  '
Set-TextValue 'D6' '581.26'
Set-TextValue 'E6' '  -1.66%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'E8' '  -0.71%  '
Set-TextValue 'D9' '3.270.74'
Set-TextValue 'E9' '  -1.21%  '
Set-TextValue 'E10' '  -4.10%  '
Set-TextValue 'E11' '  -2.21%  '
Set-TextValue 'E12' '  -3.07%  '
Set-TextValue 'D13' '3.840.04'
Set-TextValue 'E13' '  -1.04%  '
Set-TextValue 'D14' '0.138'
Set-TextValue 'E14' '  +0.04%  '
Set-TextValue 'D15' '27.56'
Set-TextValue 'E15' '  -5.47%  '
Set-TextValue 'D16' '67.952.31'
Set-TextValue 'E16' '  -1.47%  '
Set-TextValue 'E17' '  -2.96%  '
Set-TextValue 'D18' '3.227.24'
Set-TextValue 'E18' '  -2.04%  '
Set-TextValue 'D19' '5.75'
Set-TextValue 'E19' '  -2.91%  '
Set-TextValue 'D20' '13.57'
Set-TextValue 'E20' '  -1.64%  '
Set-TextValue 'D21' '398.39'
Set-TextValue 'E21' '  +2.92%  '
Set-TextValue 'D22' '7.63'
Set-TextValue 'E22' '  -2.67%  '
Set-TextValue 'B23' 'Dai'
Set-TextValue 'C23' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D23' '1.00'
Set-TextValue 'E23' '  +0.10%  '
Set-TextValue 'B24' 'Litecoin'
Set-TextValue 'C24' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D24' '71.46'
Set-TextValue 'E24' '  -0.61%  '
Set-TextValue 'D25' '0.513'
Set-TextValue 'E25' '  -1.68%  '
Set-TextValue 'E26' '  -4.67%  '
Set-TextValue 'E27' '  -1.33%  '
Set-TextValue 'E28' '  -4.21%  '
Set-TextValue 'E29' '  +0.35%  '
Set-TextValue 'E30' '  -2.81%  '
Set-TextValue 'D31' '22.70'
Set-TextValue 'E31' '  -2.12%  '
Set-TextValue 'E32' '  -6.75%  '
Set-TextValue 'D33' '6.97'
Set-TextValue 'E33' '  -4.09%  '
Set-TextValue 'E34' '  -5.97%  '
Set-TextValue 'D35' '0.998'
Set-TextValue 'E35' '  +0.00%  '
Set-TextValue 'D36' '162.78'
Set-TextValue 'E36' '  -0.58%  '
Set-TextValue 'E37' '  -5.75%  '
Set-TextValue 'E38' '  -0.47%  '
Set-TextValue 'D39' '26.82'
Set-TextValue 'E39' '  -0.64%  '
Set-TextValue 'D40' '0.811'
Set-TextValue 'E40' '  -3.77%  '
Set-TextValue 'E41' '  -2.39%  '
Set-TextValue 'D42' '6.44'
Set-TextValue 'E42' '  -4.88%  '
Set-TextValue 'D43' '2.677.81'
Set-TextValue 'E43' '  +0.73%  '
Set-TextValue 'D44' '0.0687'
Set-TextValue 'E44' '  -1.87%  '
Set-TextValue 'E45' '  -1.87%  '
Set-TextValue 'E46' '  -8.55%  '
Set-TextValue 'D47' '24.79'
Set-TextValue 'E47' '  -4.56%  '
Set-TextValue 'D48' '332.73'
Set-TextValue 'E48' '  -3.05%  '
Set-TextValue 'E49' '  -3.36%  '
Set-TextValue 'E50' '  +0.61%  '
Set-TextValue 'D51' '0.102'
Set-TextValue 'E51' '  -1.55%  '
